$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "current page (bookmark)" for "A Student's Guide to Bayesian Statistics" (row 12)
$ws.Range("C12").Value = 12

# Update selection to C20 (matches diff's new selection target)
$ws.Range("C20").Select()
